# Nina's resubmission response has now been collected - record it in the
# tracker sheet: fill in the "response_collected" file for the existing
# Sep-22-2023 row, and add a new row for the resubmission pairwise task.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (original Nina submission) - response file just came in
$ws.Range("E12").Value = "master_worker_response_tracke_Sep-22-2023.csv"

# Row 13 (new) - Nina's resubmission tracking row
$ws.Range("E13").Value = "master_worker_response_tracke_resub_Sep-22-2023.csv"
$ws.Range("A13").Value = "paiewise_resub"
$ws.Range("D13").Value = "all_submitted_tracker_nina_resubSep-22-2023.csv"
$ws.Range("F13").Value = "master_all_responses_Sep-22-2023_to_resub_Sep-22-2023_Nina.csv"

# Date column: type "Sep-22-2023" as literal text (not an auto-converted
# date) by computing it via a formula and pasting the result back as a
# value - mirrors how the rest of the Date column is stored as text.
$ws.Range("B13").Formula = '="Sep-22-2023"'
$ws.Range("B13").Copy()
$ws.Range("B13").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("C13").Value = "NV"

$ws.Range("E15").Select()
